$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 878.6
$ws.Range("I6").Value = 948.2857
$ws.Range("J6").Value = 716
$ws.Range("K6").Value = 2844.8571
$ws.Range("L6").Value = 2148
$ws.Range("M6").Value = -2732.8571
$ws.Range("N6").Value = -2372
$ws.Range("H28").Value = 4795.125
$ws.Range("I28").Value = 1349.75
$ws.Range("K28").Value = 1349.75
$ws.Range("M28").Value = -864.75
$ws.Range("H33").Value = 261.14285
$ws.Range("I33").Value = 259.67743
$ws.Range("J33").Value = 272.5
$ws.Range("K33").Value = 259.67743
$ws.Range("L33").Value = 272.5
$ws.Range("M33").Value = -30.67743000000002
$ws.Range("N33").Value = -730.5
$ws.Range("H40").Value = 125003060
$ws.Range("I40").Value = 4750
$ws.Range("J40").Value = 166669150
$ws.Range("K40").Value = 4750
$ws.Range("L40").Value = 166669150
$ws.Range("M40").Value = -4575
$ws.Range("N40").Value = -166669500
$ws.Range("H45").Value = 1888
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = $null
$ws.Range("H51").Value = 5891.2173
$ws.Range("I51").Value = 4999.5
$ws.Range("J51").Value = 5976.143
$ws.Range("K51").Value = 4999.5
$ws.Range("L51").Value = 5976.143
$ws.Range("M51").Value = -4515.5
$ws.Range("N51").Value = -6944.143
$ws.Range("H100").Value = 1744.9375
$ws.Range("I100").Value = 933.5454999999999
$ws.Range("J100").Value = 3530
$ws.Range("K100").Value = 933.5454999999999
$ws.Range("L100").Value = 3530
$ws.Range("M100").Value = -392.5454999999999
$ws.Range("N100").Value = -4612
$ws.Range("H103").Value = 20834204
$ws.Range("J103").Value = 41667670
$ws.Range("L103").Value = 125003010
$ws.Range("N103").Value = -125004182
$ws.Range("H111").Value = 10426.833
$ws.Range("I111").Value = 4126.6665
$ws.Range("K111").Value = 12379.9995
$ws.Range("M111").Value = -9312.999500000002
$ws.Range("H113").Value = 18749.25
$ws.Range("I113").Value = 24166
$ws.Range("J113").Value = 2499
$ws.Range("K113").Value = 24166
$ws.Range("L113").Value = 2499
$ws.Range("M113").Value = -20912
$ws.Range("N113").Value = -9007
$ws.Range("H132").Value = 3463.3845
$ws.Range("I132").Value = 2513.2554
$ws.Range("J132").Value = 12394.6
$ws.Range("K132").Value = 7539.7662
$ws.Range("L132").Value = 37183.8
$ws.Range("M132").Value = -5009.7662
$ws.Range("N132").Value = -42243.8
$ws.Range("H137").Value = 2582
$ws.Range("I137").Value = 1762.3334
$ws.Range("K137").Value = 5287.0002
$ws.Range("M137").Value = -2737.0002
$ws.Range("H138").Value = 5182.0713
$ws.Range("J138").Value = 19959
$ws.Range("L138").Value = 59877
$ws.Range("N138").Value = -70157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3189.9565
$ws.Range("I2").Value = 2947.0527
$ws.Range("J2").Value = 4343.75
$ws.Range("K2").Value = 2947.0527
$ws.Range("L2").Value = 4343.75
$ws.Range("M2").Value = -2834.0527
$ws.Range("N2").Value = -4569.75
$ws.Range("H5").Value = 939.1667
$ws.Range("I5").Value = 728
$ws.Range("K5").Value = 728
$ws.Range("M5").Value = -616
$ws.Range("H32").Value = 2120.077
$ws.Range("I32").Value = 2069.18
$ws.Range("J32").Value = 3392.5
$ws.Range("K32").Value = 2069.18
$ws.Range("L32").Value = 3392.5
$ws.Range("M32").Value = -1782.18
$ws.Range("N32").Value = -3966.5
$ws.Range("I40").Value = 40000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 40000
$ws.Range("L40").Value = $null
$ws.Range("M40").Value = -39824
$ws.Range("N40").Value = 0
$ws.Range("H60").Value = 63568
$ws.Range("I60").Value = 58363.57
$ws.Range("K60").Value = 58363.57
$ws.Range("M60").Value = -57630.57
$ws.Range("H63").Value = 3550
$ws.Range("I63").Value = 3600
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 3600
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -2914
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 3550
$ws.Range("I66").Value = 3600
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 18000
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -14568
$ws.Range("N66").Value = -24364
$ws.Range("H88").Value = 3186.1667
$ws.Range("I88").Value = 3281.5
$ws.Range("J88").Value = 3138.5
$ws.Range("K88").Value = 3281.5
$ws.Range("L88").Value = 3138.5
$ws.Range("M88").Value = -2875.5
$ws.Range("N88").Value = -3950.5
$ws.Range("H91").Value = 3186.1667
$ws.Range("I91").Value = 3281.5
$ws.Range("J91").Value = 3138.5
$ws.Range("K91").Value = 3281.5
$ws.Range("L91").Value = 3138.5
$ws.Range("M91").Value = -1877.5
$ws.Range("N91").Value = -5946.5
$ws.Range("H97").Value = 1585.7778
$ws.Range("I97").Value = 1409
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1409
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -913
$ws.Range("N97").Value = -3992
$ws.Range("H102").Value = 1671.65
$ws.Range("J102").Value = 2881.3333
$ws.Range("L102").Value = 2881.3333
$ws.Range("N102").Value = -6125.3333
$ws.Range("H116").Value = 3189.9565
$ws.Range("I116").Value = 2947.0527
$ws.Range("J116").Value = 4343.75
$ws.Range("K116").Value = 2947.0527
$ws.Range("L116").Value = 4343.75
$ws.Range("M116").Value = -653.0527000000002
$ws.Range("N116").Value = -8931.75
$ws.Range("H122").Value = 6059.8
$ws.Range("I122").Value = 6059.8
$ws.Range("K122").Value = 18179.4
$ws.Range("M122").Value = -15729.4
$ws.Range("H132").Value = 16669056
$ws.Range("I132").Value = 2062.4348
$ws.Range("J132").Value = 46156812
$ws.Range("K132").Value = 6187.3044
$ws.Range("L132").Value = 138470436
$ws.Range("M132").Value = -3657.3044
$ws.Range("N132").Value = -138475496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3189.9565
$ws.Range("I3").Value = 2947.0527
$ws.Range("J3").Value = 4343.75
$ws.Range("K3").Value = 2947.0527
$ws.Range("L3").Value = 4343.75
$ws.Range("M3").Value = -2833.0527
$ws.Range("N3").Value = -4571.75
$ws.Range("H4").Value = 939.1667
$ws.Range("I4").Value = 728
$ws.Range("K4").Value = 728
$ws.Range("M4").Value = -613
$ws.Range("H20").Value = 20119.5
$ws.Range("I20").Value = 27485.143
$ws.Range("J20").Value = 2933
$ws.Range("K20").Value = 27485.143
$ws.Range("L20").Value = 2933
$ws.Range("M20").Value = -27238.143
$ws.Range("N20").Value = -3427
$ws.Range("H22").Value = 1600.091
$ws.Range("I22").Value = 1690
$ws.Range("J22").Value = 1525.1666
$ws.Range("K22").Value = 1690
$ws.Range("L22").Value = 1525.1666
$ws.Range("M22").Value = -1517
$ws.Range("N22").Value = -1871.1666
$ws.Range("H86").Value = 4017.5557
$ws.Range("I86").Value = 1214.875
$ws.Range("J86").Value = 6259.7
$ws.Range("K86").Value = 1214.875
$ws.Range("L86").Value = 6259.7
$ws.Range("M86").Value = -91.875
$ws.Range("N86").Value = -8505.700000000001
$ws.Range("H89").Value = 4017.5557
$ws.Range("I89").Value = 1214.875
$ws.Range("J89").Value = 6259.7
$ws.Range("K89").Value = 6074.375
$ws.Range("L89").Value = 31298.5
$ws.Range("M89").Value = -458.375
$ws.Range("N89").Value = -42530.5
$ws.Range("H94").Value = 2902.5334
$ws.Range("I94").Value = 3093.6924
$ws.Range("K94").Value = 3093.6924
$ws.Range("M94").Value = -2642.6924
$ws.Range("H99").Value = 2394.6667
$ws.Range("I99").Value = 2196.7144
$ws.Range("J99").Value = 3087.5
$ws.Range("K99").Value = 2196.7144
$ws.Range("L99").Value = 3087.5
$ws.Range("M99").Value = -698.7143999999998
$ws.Range("N99").Value = -6083.5
$ws.Range("H107").Value = 3975.5557
$ws.Range("I107").Value = 4631.1333
$ws.Range("J107").Value = 697.6667
$ws.Range("K107").Value = 4631.1333
$ws.Range("L107").Value = 697.6667
$ws.Range("M107").Value = -2711.1333
$ws.Range("N107").Value = -4537.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 10932.667
$ws.Range("J5").Value = 10932.667
$ws.Range("L5").Value = 10932.667
$ws.Range("N5").Value = -11156.667
$ws.Range("H7").Value = 458.55554
$ws.Range("I7").Value = 454.83334
$ws.Range("J7").Value = 466
$ws.Range("K7").Value = 454.83334
$ws.Range("L7").Value = 466
$ws.Range("M7").Value = -341.83334
$ws.Range("N7").Value = -692
$ws.Range("H8").Value = 5749
$ws.Range("J8").Value = 5749
$ws.Range("L8").Value = 5749
$ws.Range("N8").Value = -6029
$ws.Range("H37").Value = 4000
$ws.Range("I37").Value = 4000
$ws.Range("K37").Value = 4000
$ws.Range("M37").Value = -3893
$ws.Range("H62").Value = 3464.3
$ws.Range("I62").Value = 3107.875
$ws.Range("K62").Value = 3107.875
$ws.Range("M62").Value = -2483.875
$ws.Range("H65").Value = 3464.3
$ws.Range("I65").Value = 3107.875
$ws.Range("K65").Value = 15539.375
$ws.Range("M65").Value = -12419.375
$ws.Range("H105").Value = 1747.5714
$ws.Range("I105").Value = 1559
$ws.Range("K105").Value = 1559
$ws.Range("M105").Value = 188
$ws.Range("H132").Value = 2935.182
$ws.Range("I132").Value = 2924.625
$ws.Range("K132").Value = 8773.875
$ws.Range("M132").Value = -6243.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1161.6
$ws.Range("I8").Value = 1161.6
$ws.Range("K8").Value = 3484.8
$ws.Range("M8").Value = -3345.8
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = $null
$ws.Range("N37").Value = 0
$ws.Range("H94").Value = 9435
$ws.Range("I94").Value = 948
$ws.Range("J94").Value = 22165.5
$ws.Range("K94").Value = 2844
$ws.Range("L94").Value = 66496.5
$ws.Range("M94").Value = -2168
$ws.Range("N94").Value = -67848.5
$ws.Range("H121").Value = 4351670
$ws.Range("I121").Value = 422
$ws.Range("K121").Value = 1266
$ws.Range("M121").Value = 44
$ws.Range("H136").Value = 6578.4443
$ws.Range("I136").Value = 1839.1428
$ws.Range("K136").Value = 5517.428400000001
$ws.Range("M136").Value = -417.4284000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1050.75
$ws.Range("I2").Value = 1064.7273
$ws.Range("J2").Value = 897
$ws.Range("K2").Value = 1064.7273
$ws.Range("L2").Value = 897
$ws.Range("M2").Value = -951.7273
$ws.Range("N2").Value = -1123
$ws.Range("H23").Value = 15349.75
$ws.Range("I23").Value = 1400
$ws.Range("J23").Value = 19999.666
$ws.Range("K23").Value = 1400
$ws.Range("L23").Value = 19999.666
$ws.Range("M23").Value = -1177
$ws.Range("N23").Value = -20445.666
$ws.Range("H70").Value = 7819.125
$ws.Range("J70").Value = 7896.4443
$ws.Range("L70").Value = 7896.4443
$ws.Range("N70").Value = -8436.444299999999
$ws.Range("H73").Value = 7819.125
$ws.Range("J73").Value = 7896.4443
$ws.Range("L73").Value = 7896.4443
$ws.Range("N73").Value = -9768.444299999999
$ws.Range("H80").Value = 3584.4
$ws.Range("I80").Value = 1494.5
$ws.Range("K80").Value = 1494.5
$ws.Range("M80").Value = -496.5
$ws.Range("H83").Value = 3584.4
$ws.Range("I83").Value = 1494.5
$ws.Range("K83").Value = 7472.5
$ws.Range("M83").Value = -2480.5
$ws.Range("H122").Value = 1969.1428
$ws.Range("I122").Value = 2418.4
$ws.Range("K122").Value = 7255.200000000001
$ws.Range("M122").Value = -4805.200000000001
$ws.Range("H126").Value = 2853.2856
$ws.Range("I126").Value = 2182.75
$ws.Range("K126").Value = 6548.25
$ws.Range("M126").Value = -4078.25
$ws.Range("H132").Value = 4002745.5
$ws.Range("I132").Value = 2574
$ws.Range("K132").Value = 7722
$ws.Range("M132").Value = -5192

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 20499.5
$ws.Range("I34").Value = 21000
$ws.Range("J34").Value = 19999
$ws.Range("K34").Value = 21000
$ws.Range("L34").Value = 19999
$ws.Range("M34").Value = -20828
$ws.Range("N34").Value = -20343
$ws.Range("H46").Value = 966.25
$ws.Range("I46").Value = 1044
$ws.Range("J46").Value = 888.5
$ws.Range("K46").Value = 1044
$ws.Range("L46").Value = 888.5
$ws.Range("M46").Value = -856
$ws.Range("N46").Value = -1264.5
$ws.Range("H93").Value = 4632295.5
$ws.Range("I93").Value = 2333
$ws.Range("J93").Value = 9262258
$ws.Range("K93").Value = 2333
$ws.Range("L93").Value = 9262258
$ws.Range("M93").Value = -1085
$ws.Range("N93").Value = -9264754
$ws.Range("H100").Value = 12516782
$ws.Range("I100").Value = 3959.9167
$ws.Range("J100").Value = 31286014
$ws.Range("K100").Value = 3959.9167
$ws.Range("L100").Value = 31286014
$ws.Range("M100").Value = -3418.9167
$ws.Range("N100").Value = -31287096
$ws.Range("H115").Value = 85099
$ws.Range("J115").Value = 85099
$ws.Range("L115").Value = 85099
$ws.Range("N115").Value = -87449
$ws.Range("H132").Value = 2121.6667
$ws.Range("I132").Value = 1836
$ws.Range("J132").Value = 2815.4285
$ws.Range("K132").Value = 5508
$ws.Range("L132").Value = 8446.2855
$ws.Range("M132").Value = -2978
$ws.Range("N132").Value = -13506.2855
$ws.Range("H136").Value = 2467.7812
$ws.Range("I136").Value = 2341.7307
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 7025.1921
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -4475.1921
$ws.Range("N136").Value = -14142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2832.5
$ws.Range("I81").Value = 2832.5
$ws.Range("K81").Value = 5665
$ws.Range("M81").Value = -4604
$ws.Range("H82").Value = 92399.5
$ws.Range("J82").Value = 92399.5
$ws.Range("L82").Value = 92399.5
$ws.Range("N82").Value = -93165.5
$ws.Range("H84").Value = 2832.5
$ws.Range("I84").Value = 2832.5
$ws.Range("K84").Value = 28325
$ws.Range("M84").Value = -23021
$ws.Range("H85").Value = 92399.5
$ws.Range("J85").Value = 92399.5
$ws.Range("L85").Value = 92399.5
$ws.Range("N85").Value = -95051.5
$ws.Range("H86").Value = 83299
$ws.Range("J86").Value = 83299
$ws.Range("L86").Value = 83299
$ws.Range("N86").Value = -85545
$ws.Range("H89").Value = 83299
$ws.Range("J89").Value = 83299
$ws.Range("L89").Value = 416495
$ws.Range("N89").Value = -427727
$ws.Range("H122").Value = 4464
$ws.Range("I122").Value = 4464
$ws.Range("K122").Value = 13392
$ws.Range("M122").Value = -10942
$ws.Range("H132").Value = 245414.7
$ws.Range("I132").Value = 1470.5151
$ws.Range("K132").Value = 4411.5453
$ws.Range("M132").Value = -1881.5453
$ws.Range("H136").Value = 188579.62
$ws.Range("I136").Value = 6941.9316
$ws.Range("K136").Value = 20825.7948
$ws.Range("M136").Value = -18275.7948
